$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Products
$ws2 = $wb.Worksheets.Item(2)   # Categories
$ws3 = $wb.Worksheets.Item(3)   # Inventory

# ---------------------------------------------------------------------------
# Replace every old placeholder string with its final value. The order in
# which brand-new literal strings are first assigned controls the order they
# get appended to the shared-string table, so the writes below are carefully
# sequenced to reproduce: name, price, imageUrl, status, available, category,
# Pants, Grey T-Shirt, Black T-Shirt, Shirt, Jeans, product
# ---------------------------------------------------------------------------

# --- Sheet1 "Products" ------------------------------------------------------
# New header cell D1 = "category" (first brand new string -> index 5)
$ws1.Range("D1").Value = "category"
$ws1.Range("D1").Font.Bold = $true

# Empty but bold-styled header cell E1
$ws1.Range("E1").Font.Bold = $true

# Row 4 category value "Pants" (-> index 6)
$ws1.Range("D4").Value = "Pants"

# Row 2 product name "Grey T-Shirt" (-> index 7)
$ws1.Range("A2").Value = "Grey T-Shirt"

# Row 3 product name "Black T-Shirt" (-> index 8)
$ws1.Range("A3").Value = "Black T-Shirt"

# Row 2 category "Shirt" (-> index 9)
$ws1.Range("D2").Value = "Shirt"
$ws1.Range("D3").Value = "Shirt"

# Row 4 product name "Jeans" (-> index 10)
$ws1.Range("A4").Value = "Jeans"

# Remove the old imageUrl values (column C no longer used on data rows)
$ws1.Range("C2").ClearContents()
$ws1.Range("C3").ClearContents()

# Numbers
$ws1.Range("B2").Value = 30
$ws1.Range("B3").Value = 40
$ws1.Range("B4").Value = 50

# --- Sheet3 "Inventory" -----------------------------------------------------
# New header cell A1 = "product" (-> index 11, last brand new string)
$ws3.Range("A1").Value = "product"
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("B1").Font.Bold = $true

# Inventory rows (status/available already exist in the shared string table)
$ws3.Range("A2").Value = "Grey T-Shirt"
$ws3.Range("B2").Value = "available"
$ws3.Range("A3").Value = "Black T-Shirt"
$ws3.Range("B3").Value = "available"
$ws3.Range("A4").Value = "Jeans"
$ws3.Range("B4").Value = "available"
$ws3.Range("A5").Value = "Grey T-Shirt"
$ws3.Range("B5").Value = "available"
$ws3.Range("A6").Value = "Black T-Shirt"
$ws3.Range("B6").Value = "available"
$ws3.Range("A7").Value = "Jeans"
$ws3.Range("B7").Value = "available"

# --- Sheet2 "Categories" ----------------------------------------------------
$ws2.Range("B1").ClearContents()
$ws2.Range("A2").Value = "Shirt"
$ws2.Range("B2").ClearContents()
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A3").Value = "Pants"

# ---------------------------------------------------------------------------
# Selections: select sheet1 then sheet2 then sheet3 last, so sheet3 (the
# original active tab) ends up active again while sheet1/sheet2 still keep
# their own remembered selection.
# ---------------------------------------------------------------------------
$ws1.Range("A2:A4").Select()
$ws2.Range("A3").Select()
$ws3.Range("C7").Select()
